$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace placeholder / joke text in the backlog with real content
$ws.Range("F3").Value  = 'Diretório para monitoramento e versionamento do projeto'
$ws.Range("F4").Value  = 'Configurar acesso dos membros da equipe'
$ws.Range("F7").Value  = 'Estabelecimento de parâmetros como missão, solução, etc.'
$ws.Range("F8").Value  = 'Desenvolvimento de logotipo adequado para a solução oferecida'
$ws.Range("F17").Value = 'Sistematizar o BackLog na plataforma Trello'
$ws.Range("F19").Value = 'Diagramar a solução oferecida, destacando a importância do sensor e do monitoramento'
$ws.Range("F20").Value = 'Plataforma de trabalho essencial para o desenvolvimento do Banco de Dados'
$ws.Range("E24").Value = 'Scripts apropriados para inserir dados'
$ws.Range("F24").Value = 'Elaborar linhas de código em MySQL para inserir dados nas tabelas adequadas.'
$ws.Range("E25").Value = 'Scripts apropriados para retornar dados'
$ws.Range("F25").Value = 'Elaborar linhas de código em MySQL para selecionar dados das tabelas adequadas.'

# Adjust the sheet view: zoom to 80% and move the active selection
$ws.Application.ActiveWindow.Zoom = 80
$ws.Range("L8").Select()

# Auto-fit column F (and H) to account for the new, longer text
$ws.Columns("F:F").AutoFit() | Out-Null
$ws.Columns("H:H").AutoFit() | Out-Null

$wb.Save()
